$d = $word.ActiveDocument

# The "Schedule" table has a "Due" column (column 3). Several rows have an
# empty paragraph in that cell that is missing the "Compact" paragraph
# style used everywhere else in the table. Apply the "Compact" style to
# every empty paragraph in that column so it matches the rest of the table.
$table = $d.Tables.Item(1)

for ($r = 1; $r -le $table.Rows.Count; $r++) {
    $cell = $table.Cell($r, 3)
    foreach ($p in $cell.Range.Paragraphs) {
        # Paragraph text includes trailing paragraph/cell marks (chr 13 / 7),
        # so strip those control characters before checking for emptiness.
        $clean = $p.Range.Text -replace "[\x07\r\n]", ""
        if ($clean -eq "") {
            $p.Style = "Compact"
        }
    }
}
